$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$theDate = Get-Date -Year 2014 -Month 10 -Day 15 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

# Fix the date typo in row 21 (A21): 2015-10-14 -> 2014-10-15
$ws.Range("A21").Value = $theDate

# Fill in row 22: preparing for presentation
$ws.Range("A22").Value = $theDate
$ws.Range("B22").Value = 0.5
$ws.Range("C22").Value = 0.58333333333333337
$ws.Range("E22").Value = "preparing for presentation"

# Fill in row 23: conducting presentation
$ws.Range("A23").Value = $theDate
$ws.Range("B23").Value = 0.63888888888888895
$ws.Range("C23").Value = 0.66666666666666663
$ws.Range("E23").Value = "conducting presentation"

# Update the active selection to E24
$ws.Range("E24").Select()
